$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 12 of the timesheet table with a new entry
$ws.Range("A12").Value = "Wk[12] Saturday 2.6.18"
$ws.Range("B12").Value = "1600 - 2300"
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = "Improving readability of code. It was so bad."

# Update the selected cell to B13
$ws.Range("B13").Select()

$wb.Save()
